# Oneida spreadsheet edit script
# Implements: richer run formatting on B1 ("Ukwehok&ha" -> adds SUPER/sub demo runs
# with bold/italic/strike/underline/superscript/subscript), font-size tweaks on a
# couple of underline runs, new demonstration rows 20/21/23/25 and a new "Sheet2"
# tab that duplicates the combiner test string with extra superscript+underline
# combo ("Fix combiner in Ahom").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

### ---- Row 1 (B1): turn the single-run "Ukwehok&" + "ha" string into a full ----
### ---- showcase of bold / italic / strikethrough / underline / superscript / ----
### ---- subscript runs, appending " SUPER sub" to the original text.         ----
$b1 = $ws1.Range("B1")
$b1.Value = "Ukwehok&ha SUPER sub"
$b1.Characters(1,2).Font.Bold = $true          # "Uk"
$b1.Characters(4,2).Font.Italic = $true        # "eh"
$b1.Characters(6,1).Font.Strikethrough = $true # "o"
$b1.Characters(9,2).Font.Underline = $true     # "ha"
$b1.Characters(12,5).Font.Superscript = $true  # "SUPER"
$b1.Characters(18,3).Font.Subscript = $true    # "sub"
$ws1.Rows.Item(1).RowHeight = 28.5

### ---- Rows 3 & 4: bump the underlined run's font size from 16 to 20 ----
$b3 = $ws1.Range("B3")
$b3.Characters(14,2).Font.Size = 20   # "ha"

$b4 = $ws1.Range("B4")
$b4.Characters(14,1).Font.Size = 20   # "e"

### ---- Row 20: shrink font to 18pt, adjust row height ----
$ws1.Range("B20").Font.Size = 18
$ws1.Rows.Item(20).RowHeight = 40.5

### ---- Row 21: switch A21 to Times New Roman 20pt, B21 to Oneida 32pt, taller row ----
$a21 = $ws1.Range("A21")
$a21.Font.Name = "Times New Roman"
$a21.Font.Size = 20

$ws1.Range("B21").Font.Size = 32
$ws1.Rows.Item(21).RowHeight = 72

### ---- Row 23 (new): duplicate of row 18's B cell, plus an annotation in C ----
$ws1.Range("B18").Copy($ws1.Range("B23"))
$ws1.Range("B23").WrapText = $true
$ws1.Rows.Item(23).RowHeight = 24.45

$c23 = $ws1.Range("C23")
$c23.Value = "Duplicate of line 18"
$c23.Font.Name = "Arial"
$c23.Font.Size = 14

### ---- Row 25 (new): duplicate of row 21 (A + B), with underline added to "su" ----
$a25 = $ws1.Range("A25")
$a25.Value = 21
$a25.Font.Name = "Times New Roman"
$a25.Font.Size = 20
$a25.Font.Underline = $true

$ws1.Range("B21").Copy($ws1.Range("B25"))
$ws1.Range("B25").WrapText = $true
$ws1.Range("B25").Font.Size = 32
$ws1.Range("B25").Characters(12,2).Font.Underline = $true   # "su"

$ws1.Rows.Item(25).RowHeight = 72

$c25 = $ws1.Range("C25")
$c25.Value = "Copy of line 21, but with underline added."

### ---- New "Sheet2" tab, placed right after Sheet1 ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Columns.Item(1).ColumnWidth = 118.5
$ws2.Rows.Item(1).RowHeight = 50.25

# Seed A1 with the Oneida 20pt base style (copied from a plain Oneida cell)
# before overwriting its value, so the font table reuses family="0" like the
# rest of the Oneida runs instead of inheriting Arial's family classification.
$ws1.Range("B2").Copy($ws2.Range("A1"))

$a1b = $ws2.Range("A1")
$a1b.Value = "Sheet2 Ukwehok&ha SUPER>^ sub<>"
$a1b.Characters(1,9).Font.Bold = $true           # "Sheet2 Uk"
$a1b.Characters(11,2).Font.Italic = $true        # "eh"
$a1b.Characters(13,1).Font.Strikethrough = $true # "o"
$a1b.Characters(16,2).Font.Underline = $true     # "ha"
$a1b.Characters(19,3).Font.Superscript = $true   # "SUP"
$a1b.Characters(22,4).Font.Superscript = $true   # "ER>^"
$a1b.Characters(22,4).Font.Underline = $true
$a1b.Characters(27,5).Font.Subscript = $true     # "sub<>"

### ---- Restore Sheet1 as the active / selected tab ----
$ws1.Activate()
$ws1.Range("B26").Select()

Write-Host "edit complete"
